# Commit: "break out stock.yaml completed"
#
# 1) Append a freshly-scraped batch (23/12/2024 11:35:04) of 21 stocks to the
#    "day" sheet, rows 1189-1209 (dimension grows from A1:I1188 to A1:I1209).
#    The bsecode column (D) is written as TEXT for this batch (matches the
#    source stock.yaml import, which kept codes as strings) except for
#    ANGELONE (row 1191), which has no BSE code and is left blank.
# 2) Normalise the "week" sheet's existing bsecode column (D674:D704) from
#    text to numeric, matching the rest of the workbook's bsecode columns.

$wb = $excel.ActiveWorkbook
$wsDay = $wb.Worksheets.Item("day")
$wsWeek = $wb.Worksheets.Item("week")

$dayRows = @(
    @(1, "ABBOTINDIA", "Abbott India Limited", "500488", -0.26, 28581, 7167, "day", "23/12/2024 11:35:04"),
    @(2, "INDIGO", "Interglobe Aviation Limited", "539448", 1.01, 4439.95, 552061, "day", "23/12/2024 11:35:04"),
    @(3, "ANGELONE", "Angel One Ltd", "", -0.77, 2859.55, 1167690, "day", "23/12/2024 11:35:04"),
    @(4, "BALKRISIND", "Balkrishna Industries Limited", "502355", 0.96, 2816.75, 208021, "day", "23/12/2024 11:35:04"),
    @(5, "POLICYBZR", "PB Fintech Ltd", "543390", 1.16, 2116.4, 1102384, "day", "23/12/2024 11:35:04"),
    @(6, "COROMANDEL", "Coromandel International Limited", "506395", 1.27, 1862.7, 881011, "day", "23/12/2024 11:35:04"),
    @(7, "HDFCBANK", "Hdfc Bank Limited", "500180", 1.67, 1801, 5522413, "day", "23/12/2024 11:35:04"),
    @(8, "BATAINDIA", "Bata India Limited", "500043", -0.61, 1334.1, 83487, "day", "23/12/2024 11:35:04"),
    @(9, "ZYDUSLIFE", "Zydus Lifesciences Ltd", "532321", -0.02, 973.35, 758156, "day", "23/12/2024 11:35:04"),
    @(10, "LICI", "Life Insurance Corporation of India", "543526", 0.34, 904.8, 1947520, "day", "23/12/2024 11:35:04"),
    @(11, "CONCOR", "Container Corporation Of India Limited", "531344", 0.92, 775, 506816, "day", "23/12/2024 11:35:04"),
    @(12, "CGPOWER", "CG Power and Industrial Solutions Ltd", "500093", -0.34, 727.6, 1507460, "day", "23/12/2024 11:35:04"),
    @(13, "TATAMOTORS", "Tata Motors Limited", "500570", -0.26, 722.2, 9653869, "day", "23/12/2024 11:35:04"),
    @(14, "APOLLOTYRE", "Apollo Tyres Limited", "500877", 0.03, 532.1, 638146, "day", "23/12/2024 11:35:04"),
    @(15, "TATAPOWER", "Tata Power Company Limited", "500400", -0.3, 399.9, 9052201, "day", "23/12/2024 11:35:04"),
    @(16, "ZOMATO", "Zomato Ltd", "543320", -2.89, 273.95, 74588008, "day", "23/12/2024 11:35:04"),
    @(17, "CUB", "City Union Bank Limited", "532210", -0.5, 174.55, 3144525, "day", "23/12/2024 11:35:04"),
    @(18, "SJVN", "Sjvn Limited", "533206", -1.59, 109.77, 4756043, "day", "23/12/2024 11:35:04"),
    @(19, "NHPC", "Nhpc Limited", "533098", 0.6899999999999999, 82.03, 14354382, "day", "23/12/2024 11:35:04"),
    @(20, "GMRAIRPORT", "GMR Airports Ltd", "532754", 0.61, 79.3, 19622277, "day", "23/12/2024 11:35:04"),
    @(21, "YESBANK", "Yes Bank Limited", "532648", 0.4, 19.91, 52739999, "day", "23/12/2024 11:35:04")
)
$weekCodes = @(
    @(674, 500530),
    @(675, 500387),
    @(676, 532538),
    @(677, 500251),
    @(678, 541154),
    @(679, 500420),
    @(680, 533150),
    @(681, 502355),
    @(682, 506401),
    @(683, 500300),
    @(684, 532343),
    @(685, 533398),
    @(686, 533309),
    @(687, 533274),
    @(688, 500247),
    @(689, 517354),
    @(690, 511243),
    @(691, 532321),
    @(692, 532286),
    @(693, 532868),
    @(694, 500093),
    @(695, 533148),
    @(696, 500670),
    @(697, 500877),
    @(698, 500295),
    @(699, 500400),
    @(700, 539876),
    @(701, 534816),
    @(702, 20712),
    @(703, 533519),
    @(704, 532648)
)
# --- Part 1: append the new batch of rows to the "day" sheet ---
$startRow = 1189
for ($i = 0; $i -lt $dayRows.Count; $i++) {
    $item = $dayRows[$i]
    $r = $startRow + $i

    $wsDay.Cells.Item($r, 1).Value = $item[0]          # sr
    $wsDay.Cells.Item($r, 2).Value = $item[1]           # nsecode
    $wsDay.Cells.Item($r, 3).Value = $item[2]           # name

    $bsecode = $item[3]
    if ($bsecode -ne "") {
        # Leading apostrophe forces text storage (mirrors the source data,
        # where bsecode came through as a string rather than a number).
        $wsDay.Cells.Item($r, 4).Formula = "'" + $bsecode
    }
    # else: no BSE code available for this stock -> leave D blank

    $wsDay.Cells.Item($r, 5).Value = $item[4]           # per_chg
    $wsDay.Cells.Item($r, 6).Value = $item[5]           # close
    $wsDay.Cells.Item($r, 7).Value = $item[6]           # volume
    $wsDay.Cells.Item($r, 8).Value = $item[7]           # timeframe
    $wsDay.Cells.Item($r, 9).Value = $item[8]           # Date Time
}

# --- Part 2: convert "week" sheet's D674:D704 bsecode from text to numeric ---
for ($i = 0; $i -lt $weekCodes.Count; $i++) {
    $pair = $weekCodes[$i]
    $r = $pair[0]
    $code = $pair[1]
    $wsWeek.Cells.Item($r, 4).Value = $code
}
